$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 4) mirroring the existing entries
$ws.Range("A4").Value = "2021-present"
$ws.Range("C4").Value = "Faculty of Science"
$ws.Range("B4").Value = "Statistician - Animal Welfare Ethical Review Body"
$ws.Range("D4").Value = "University of East Anglia"

# Move the active selection to reflect post-edit cursor position
$ws.Range("B11").Select()
